$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("Input")
$wsCases = $wb.Worksheets.Item("Cases")

# ---------------------------------------------------------------------------
# Input sheet: update the scalar test parameters in row 2. Every formula on
# the Cases sheet that reads Input!E2/F2/G2/H2/I2 recalculates automatically.
# ---------------------------------------------------------------------------
$wsInput.Range("B2").Value = 37.568
$wsInput.Range("C2").Value = 37.568
$wsInput.Range("E2").Value = 10.520333333333333
$wsInput.Range("F2").Value = 2.7475
$wsInput.Range("G2").Value = 4
$wsInput.Range("H2").Value = 5
$wsInput.Range("I2").Value = 6

# Project location name: "Fredericia" -> "Aadum"
$wsInput.Range("J2").Value = "Aadum"

# ---------------------------------------------------------------------------
# Cases sheet: clear a stray leftover formula in AQ25 (keep the cell's style)
# ---------------------------------------------------------------------------
$wsCases.Range("AQ25").ClearContents()

# ---------------------------------------------------------------------------
# Cases sheet: append new test-case row 34
# ---------------------------------------------------------------------------
$wsCases.Range("A34").Value = "3.9.9.9"
$wsCases.Range("B34").Value = "PF_test"
$wsCases.Range("C34").Value = 33
$wsCases.Range("D34").Value = 1
$wsCases.Range("E34").Value = 1

$wsCases.Range("F30").Copy()
$wsCases.Range("F34").PasteSpecial(-4122)
$wsCases.Range("F34").Formula = "=Input!I2"

$wsCases.Range("G34").Value = 1
$wsCases.Range("H34").Value = 0
$wsCases.Range("I34").Value = 1
$wsCases.Range("J34").Value = 0
$wsCases.Range("K34").Value = 0
$wsCases.Range("L34").Value = 0
$wsCases.Range("M34").Value = 0
$wsCases.Range("N34").Value = 0
$wsCases.Range("O34").Value = 0
$wsCases.Range("P34").Value = 1
$wsCases.Range("Q34").Value = 20
$wsCases.Range("R34").Value = 3
$wsCases.Range("S34").Value = 0.9
$wsCases.Range("T34").Value = 0
$wsCases.Range("U34").Value = 10
$wsCases.Range("V34").Value = -0.9
$wsCases.Range("W34").Value = 0
$wsCases.Range("X34").Value = 0
$wsCases.Range("Y34").Value = 0
$wsCases.Range("Z34").Value = 0
$wsCases.Range("AA34").Value = 0
$wsCases.Range("AB34").Value = 0
$wsCases.Range("AC34").Value = 0
$wsCases.Range("AD34").Value = 0
$wsCases.Range("AE34").Value = 0
$wsCases.Range("AF34").Value = 0
$wsCases.Range("AG34").Value = 0
$wsCases.Range("AH34").Value = 0
$wsCases.Range("AI34").Value = 0
$wsCases.Range("AJ34").Value = 0
$wsCases.Range("AK34").Value = 0
$wsCases.Range("AL34").Value = 0

# ---------------------------------------------------------------------------
# View state: active sheet/selection swaps from Input -> Cases
# ---------------------------------------------------------------------------
$wsInput.Range("O13").Select()
$wsCases.Select()
$wsCases.Range("AB25").Select()
